# Append a new daily data row (2025/10/04) to the tracking sheet,
# matching the pattern of the existing rows (row 58/59).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 60

# Column A holds the date as literal text (e.g. "2025/09/22"), not a
# real Excel date value, in every existing row. Force text formatting
# before assigning the value so Excel doesn't auto-convert the
# "yyyy/mm/dd" string into a date serial number, then reset the cell
# style back to Normal so it matches the unstyled data cells above it.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/04"
$cellA.Style = "Normal"

# Column B: day-of-week label (plain text, no special formatting needed)
$ws.Cells.Item($row, 2).Value = "土"

# Column C: time/hour value
$ws.Cells.Item($row, 3).Value = 16

# Column D: ranking value
$ws.Cells.Item($row, 4).Value = 201
